$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: status text for the 974ac0c5-... file changes from
#     "Ready for handoff" to "Handback transform failed" (shows up in both
#     E3 and F3, since the report writes the same status string into both
#     locale columns). ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- The same status string is repeated on each per-locale sheet's
#     "Status" column (column C) for that file's row, so it needs the same
#     update there. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus

# --- zh-cn sheet: the handback process recorded an error detail for the
#     second row (974ac0c5-...), because the handback file name didn't
#     match the handoff file name. Also widen the "Error Detail" column
#     (column P, the 16th column) to fit the new text. Other columns with
#     a saved width of 40 report 39.17 via ColumnWidth in this workbook
#     (width<->ColumnWidth conversion has a constant ~0.83 character
#     offset), so use that figure to land on a saved width of exactly 40. ---
$wsZhCn.Range("P3").Value = "Handback file name: 1wzx1cd1.0dq is different with handoff file name: 974ac0c5-0375-45af-8595-83a27a765dae.1a1c288cc9257d6c84343141309454a6d06536f3.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet: same error-detail report for the de-de locale, plus the
#     same column widening. ---
$wsDeDe.Range("P3").Value = "Handback file name: 1wzx1cd1.0dq is different with handoff file name: 974ac0c5-0375-45af-8595-83a27a765dae.1a1c288cc9257d6c84343141309454a6d06536f3.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
